$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4524.1724
$ws.Range("I62").Value = 3800.6667
$ws.Range("J62").Value = 5034.8823
$ws.Range("K62").Value = 3800.6667
$ws.Range("L62").Value = 5034.8823
$ws.Range("M62").Value = -3176.6667
$ws.Range("N62").Value = -6282.8823
$ws.Range("H65").Value = 4524.1724
$ws.Range("I65").Value = 3800.6667
$ws.Range("J65").Value = 5034.8823
$ws.Range("K65").Value = 19003.3335
$ws.Range("L65").Value = 25174.4115
$ws.Range("M65").Value = -15883.3335
$ws.Range("N65").Value = -31414.4115
$ws.Range("H113").Value = 4925.8423
$ws.Range("I113").Value = 3981.3635
$ws.Range("J113").Value = 6224.5
$ws.Range("K113").Value = 3981.3635
$ws.Range("L113").Value = 6224.5
$ws.Range("M113").Value = -727.3634999999999
$ws.Range("N113").Value = -12732.5
$ws.Range("H116").Value = 5786.4287
$ws.Range("I116").Value = 2362.875
$ws.Range("J116").Value = 10351.167
$ws.Range("K116").Value = 2362.875
$ws.Range("L116").Value = 10351.167
$ws.Range("M116").Value = 1079.125
$ws.Range("N116").Value = -17235.167
$ws.Range("H132").Value = 4858.4443
$ws.Range("I132").Value = 2272.3447
$ws.Range("J132").Value = 15572.286
$ws.Range("K132").Value = 6817.034100000001
$ws.Range("L132").Value = 46716.858
$ws.Range("M132").Value = -4287.034100000001
$ws.Range("N132").Value = -51776.858
$ws.Range("H135").Value = 25000764
$ws.Range("I135").Value = 565.0769
$ws.Range("J135").Value = 71429704
$ws.Range("K135").Value = 5085.6921
$ws.Range("L135").Value = 642867336
$ws.Range("M135").Value = -2550.6921
$ws.Range("N135").Value = -642872406
$ws.Range("H137").Value = 2107.3684
$ws.Range("I137").Value = 1559
$ws.Range("J137").Value = 2716.6667
$ws.Range("K137").Value = 4677
$ws.Range("L137").Value = 8150.000100000001
$ws.Range("M137").Value = -2127
$ws.Range("N137").Value = -13250.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1582.25
$ws.Range("I61").Value = 1484.8462
$ws.Range("J61").Value = 1763.1428
$ws.Range("K61").Value = 1484.8462
$ws.Range("L61").Value = 1763.1428
$ws.Range("M61").Value = -1272.8462
$ws.Range("N61").Value = -2187.1428
$ws.Range("H74").Value = 24725.861
$ws.Range("I74").Value = 41080.64
$ws.Range("J74").Value = 2010.8889
$ws.Range("K74").Value = 41080.64
$ws.Range("L74").Value = 2010.8889
$ws.Range("M74").Value = -40206.64
$ws.Range("N74").Value = -3758.8889
$ws.Range("H77").Value = 24725.861
$ws.Range("I77").Value = 41080.64
$ws.Range("J77").Value = 2010.8889
$ws.Range("K77").Value = 205403.2
$ws.Range("L77").Value = 10054.4445
$ws.Range("M77").Value = -201035.2
$ws.Range("N77").Value = -18790.4445
$ws.Range("H132").Value = 1752.5319
$ws.Range("I132").Value = 1716.4634
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 5149.3902
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -2619.3902
$ws.Range("N132").Value = -11057
$ws.Range("H136").Value = 1582.25
$ws.Range("I136").Value = 1484.8462
$ws.Range("J136").Value = 1763.1428
$ws.Range("K136").Value = 4454.5386
$ws.Range("L136").Value = 5289.428400000001
$ws.Range("M136").Value = -1904.5386
$ws.Range("N136").Value = -10389.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1906.1765
$ws.Range("I99").Value = 1899.4546
$ws.Range("J99").Value = 1918.5
$ws.Range("K99").Value = 1899.4546
$ws.Range("L99").Value = 1918.5
$ws.Range("M99").Value = -401.4546
$ws.Range("N99").Value = -4914.5
$ws.Range("H105").Value = 1901.0723
$ws.Range("I105").Value = 1639.2075
$ws.Range("K105").Value = 1639.2075
$ws.Range("M105").Value = 107.7925
$ws.Range("H134").Value = 628145.4399999999
$ws.Range("I134").Value = 1114463.9
$ws.Range("J134").Value = 2878.8928
$ws.Range("K134").Value = 3343391.7
$ws.Range("L134").Value = 8636.678400000001
$ws.Range("M134").Value = -3340856.7
$ws.Range("N134").Value = -13706.6784

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 909.625
$ws.Range("I16").Value = 735.6
$ws.Range("J16").Value = 1199.6666
$ws.Range("K16").Value = 735.6
$ws.Range("L16").Value = 1199.6666
$ws.Range("M16").Value = -448.6
$ws.Range("N16").Value = -1773.6666
$ws.Range("H31").Value = 1377.8209
$ws.Range("I31").Value = 869.5
$ws.Range("J31").Value = 2491.2856
$ws.Range("K31").Value = 869.5
$ws.Range("L31").Value = 2491.2856
$ws.Range("M31").Value = -574.5
$ws.Range("N31").Value = -3081.2856
$ws.Range("H34").Value = 1377.8209
$ws.Range("I34").Value = 869.5
$ws.Range("J34").Value = 2491.2856
$ws.Range("K34").Value = 869.5
$ws.Range("L34").Value = 2491.2856
$ws.Range("M34").Value = -667.5
$ws.Range("N34").Value = -2895.2856
$ws.Range("H94").Value = 1039.9333
$ws.Range("I94").Value = 1179.8
$ws.Range("J94").Value = 970
$ws.Range("K94").Value = 1179.8
$ws.Range("L94").Value = 970
$ws.Range("M94").Value = -728.8
$ws.Range("N94").Value = -1872
$ws.Range("H99").Value = 3681.923
$ws.Range("I99").Value = 3542.4443
$ws.Range("J99").Value = 3995.75
$ws.Range("K99").Value = 3542.4443
$ws.Range("L99").Value = 3995.75
$ws.Range("M99").Value = -2044.4443
$ws.Range("N99").Value = -6991.75
$ws.Range("H113").Value = 909.625
$ws.Range("I113").Value = 735.6
$ws.Range("J113").Value = 1199.6666
$ws.Range("K113").Value = 735.6
$ws.Range("L113").Value = 1199.6666
$ws.Range("M113").Value = 1434.4
$ws.Range("N113").Value = -5539.6666
$ws.Range("H126").Value = 3681.923
$ws.Range("I126").Value = 3542.4443
$ws.Range("J126").Value = 3995.75
$ws.Range("K126").Value = 10627.3329
$ws.Range("L126").Value = 11987.25
$ws.Range("M126").Value = -8157.332900000001
$ws.Range("N126").Value = -16927.25
$ws.Range("H132").Value = 1030800.9
$ws.Range("I132").Value = 2079.1724
$ws.Range("J132").Value = 5292648
$ws.Range("K132").Value = 6237.5172
$ws.Range("L132").Value = 15877944
$ws.Range("M132").Value = -3707.5172
$ws.Range("N132").Value = -15883004
$ws.Range("H134").Value = 2086.5
$ws.Range("I134").Value = 2205.1936
$ws.Range("J134").Value = 1677.6666
$ws.Range("K134").Value = 6615.5808
$ws.Range("L134").Value = 5032.9998
$ws.Range("M134").Value = -4080.5808
$ws.Range("N134").Value = -10102.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2525768.2
$ws.Range("I113").Value = 3788402.5
$ws.Range("J113").Value = 499.25
$ws.Range("K113").Value = 11365207.5
$ws.Range("L113").Value = 1497.75
$ws.Range("M113").Value = -11363037.5
$ws.Range("N113").Value = -5837.75
$ws.Range("H122").Value = 581.5714
$ws.Range("I122").Value = 540.8
$ws.Range("J122").Value = 683.5
$ws.Range("K122").Value = 4867.2
$ws.Range("L122").Value = 6151.5
$ws.Range("M122").Value = -2417.2
$ws.Range("N122").Value = -11051.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 614.7778
$ws.Range("I97").Value = 633.5
$ws.Range("J97").Value = 532.4
$ws.Range("K97").Value = 633.5
$ws.Range("L97").Value = 532.4
$ws.Range("M97").Value = -137.5
$ws.Range("N97").Value = -1524.4
$ws.Range("H132").Value = 2176600.5
$ws.Range("I132").Value = 2692.4688
$ws.Range("J132").Value = 7145533.5
$ws.Range("K132").Value = 8077.4064
$ws.Range("L132").Value = 21436600.5
$ws.Range("M132").Value = -5547.4064
$ws.Range("N132").Value = -21441660.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 860.0909
$ws.Range("I22").Value = 733.3333
$ws.Range("J22").Value = 907.625
$ws.Range("K22").Value = 733.3333
$ws.Range("L22").Value = 907.625
$ws.Range("M22").Value = -438.3333
$ws.Range("N22").Value = -1497.625
$ws.Range("H27").Value = 860.0909
$ws.Range("I27").Value = 733.3333
$ws.Range("J27").Value = 907.625
$ws.Range("K27").Value = 733.3333
$ws.Range("L27").Value = 907.625
$ws.Range("M27").Value = -626.3333
$ws.Range("N27").Value = -1121.625
$ws.Range("H136").Value = 1427.7273
$ws.Range("I136").Value = 1004.5
$ws.Range("J136").Value = 2168.375
$ws.Range("K136").Value = 3013.5
$ws.Range("L136").Value = 6505.125
$ws.Range("M136").Value = -463.5
$ws.Range("N136").Value = -11605.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1143.8833
$ws.Range("I136").Value = 669
$ws.Range("J136").Value = 2345.0588
$ws.Range("K136").Value = 2007
$ws.Range("L136").Value = 7035.176399999999
$ws.Range("M136").Value = 543
$ws.Range("N136").Value = -12135.1764
